$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 12.79505811484147
$ws.Cells.Item(2, 3).Value = 9.253429113418647
$ws.Cells.Item(2, 4).Value = 4.057763500715986
$ws.Cells.Item(2, 5).Value = 11.54444127637291
$ws.Cells.Item(2, 6).Value = 21.21162353542267
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 13).Value = 14.5377232406679
$ws.Cells.Item(2, 15).Value = 18.84733671906806
$ws.Cells.Item(3, 2).Value = 12.14100485638195
$ws.Cells.Item(3, 3).Value = 8.82313187892489
$ws.Cells.Item(3, 4).Value = 4.015209686612629
$ws.Cells.Item(3, 5).Value = 11.47080857032227
$ws.Cells.Item(3, 6).Value = 21.20815380122425
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 13).Value = 14.22207651298682
$ws.Cells.Item(3, 15).Value = 18.91520893882586
$ws.Cells.Item(4, 2).Value = 11.72132454061443
$ws.Cells.Item(4, 3).Value = 8.54667729497978
$ws.Cells.Item(4, 4).Value = 3.988735954430699
$ws.Cells.Item(4, 5).Value = 11.43049684559042
$ws.Cells.Item(4, 6).Value = 21.21499624123008
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 13).Value = 14.02759131164497
$ws.Cells.Item(4, 15).Value = 18.96380518993152
$ws.Cells.Item(5, 2).Value = 11.54593226385338
$ws.Cells.Item(5, 3).Value = 8.431038499509045
$ws.Cells.Item(5, 4).Value = 3.977867639515519
$ws.Cells.Item(5, 5).Value = 11.41531406874007
$ws.Cells.Item(5, 6).Value = 21.22003410938339
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 13).Value = 13.948280266646
$ws.Cells.Item(5, 15).Value = 18.98533822511592
$ws.Cells.Item(6, 2).Value = 11.51655050331577
$ws.Cells.Item(6, 3).Value = 8.41165992389779
$ws.Cells.Item(6, 4).Value = 3.976058329993128
$ws.Cells.Item(6, 5).Value = 11.41286846491397
$ws.Cells.Item(6, 6).Value = 21.22100626224119
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 13).Value = 13.93511067659381
$ws.Cells.Item(6, 15).Value = 18.98901790345012
$ws.Cells.Item(7, 2).Value = 11.71897656887415
$ws.Cells.Item(7, 3).Value = 8.545129677889578
$ws.Cells.Item(7, 4).Value = 3.988589695610097
$ws.Cells.Item(7, 5).Value = 11.43028703146089
$ws.Cells.Item(7, 6).Value = 21.21505508632099
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 13).Value = 14.02652177273792
$ws.Cells.Item(7, 15).Value = 18.96408860366121
$ws.Cells.Item(8, 2).Value = 12.5733930474712
$ws.Cells.Item(8, 3).Value = 9.107659737879139
$ws.Cells.Item(8, 4).Value = 4.043166728446545
$ws.Cells.Item(8, 5).Value = 11.51804437614503
$ws.Cells.Item(8, 6).Value = 21.20856152742332
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 13).Value = 14.42909576005054
$ws.Cells.Item(8, 15).Value = 18.86929584529765
$ws.Cells.Item(9, 2).Value = 14.09897259975212
$ws.Cells.Item(9, 3).Value = 10.10999521382413
$ws.Cells.Item(9, 4).Value = 4.147156059211882
$ws.Cells.Item(9, 5).Value = 11.72831648210886
$ws.Cells.Item(9, 6).Value = 21.26720711000982
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 13).Value = 15.20837179561112
$ws.Cells.Item(9, 15).Value = 18.73880019948323
$ws.Cells.Item(10, 2).Value = 15.12158952834518
$ws.Cells.Item(10, 3).Value = 10.78109782053133
$ws.Cells.Item(10, 4).Value = 4.221320778085648
$ws.Cells.Item(10, 5).Value = 11.90497093129739
$ws.Cells.Item(10, 6).Value = 21.35388469130671
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 13).Value = 15.76859440031214
$ws.Cells.Item(10, 15).Value = 18.67728926618764
$ws.Cells.Item(11, 2).Value = 15.56444518437461
$ws.Cells.Item(11, 3).Value = 11.0716211639669
$ws.Cells.Item(11, 4).Value = 4.254495616632545
$ws.Cells.Item(11, 5).Value = 11.98986400231102
$ws.Cells.Item(11, 6).Value = 21.40274170417204
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 13).Value = 16.01959224878216
$ws.Cells.Item(11, 15).Value = 18.65688088749958
$ws.Cells.Item(12, 2).Value = 15.72886609835303
$ws.Cells.Item(12, 3).Value = 11.17947343652028
$ws.Cells.Item(12, 4).Value = 4.266970654267523
$ws.Cells.Item(12, 5).Value = 12.02263562002417
$ws.Cells.Item(12, 6).Value = 21.42259045991287
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 13).Value = 16.113993805704
$ws.Cells.Item(12, 15).Value = 18.65024980406709
$ws.Cells.Item(13, 2).Value = 15.69360186619621
$ws.Cells.Item(13, 3).Value = 11.15634225287187
$ws.Cells.Item(13, 4).Value = 4.264287928499031
$ws.Cells.Item(13, 5).Value = 12.01555034473366
$ws.Cells.Item(13, 6).Value = 21.41825589230727
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 13).Value = 16.09369286934599
$ws.Cells.Item(13, 15).Value = 18.65162901161904
$ws.Cells.Item(14, 2).Value = 15.5780382415564
$ws.Cells.Item(14, 3).Value = 11.08053777710703
$ws.Cells.Item(14, 4).Value = 4.255523737236651
$ws.Cells.Item(14, 5).Value = 11.99254779010707
$ws.Cells.Item(14, 6).Value = 21.40434771901505
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 13).Value = 16.0273721722386
$ws.Cells.Item(14, 15).Value = 18.65631331285565
$ws.Cells.Item(15, 2).Value = 15.50682341952723
$ws.Cells.Item(15, 3).Value = 11.03382266644361
$ws.Cells.Item(15, 4).Value = 4.250143822746086
$ws.Cells.Item(15, 5).Value = 11.97853854586334
$ws.Cells.Item(15, 6).Value = 21.39600378582154
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 13).Value = 15.98666208653544
$ws.Cells.Item(15, 15).Value = 18.65932568678949
$ws.Cells.Item(16, 2).Value = 15.09219235826531
$ws.Cells.Item(16, 3).Value = 10.76181078693405
$ws.Cells.Item(16, 4).Value = 4.219140826797372
$ws.Cells.Item(16, 5).Value = 11.89951187626071
$ws.Cells.Item(16, 6).Value = 21.35088085006247
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 13).Value = 15.75210559929974
$ws.Cells.Item(16, 15).Value = 18.67877610687573
$ws.Cells.Item(17, 2).Value = 14.83205648302686
$ws.Cells.Item(17, 3).Value = 10.5911282008131
$ws.Cells.Item(17, 4).Value = 4.199972533661757
$ws.Cells.Item(17, 5).Value = 11.85217351612617
$ws.Cells.Item(17, 6).Value = 21.325609088288
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 13).Value = 15.60715907849708
$ws.Cells.Item(17, 15).Value = 18.69265403353893
$ws.Cells.Item(18, 2).Value = 14.68033370125042
$ws.Cells.Item(18, 3).Value = 10.49156821245601
$ws.Cells.Item(18, 4).Value = 4.18889482260318
$ws.Cells.Item(18, 5).Value = 11.82537407435997
$ws.Cells.Item(18, 6).Value = 21.31196139605511
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 13).Value = 15.52343239917068
$ws.Cells.Item(18, 15).Value = 18.70134864789445
$ws.Cells.Item(19, 2).Value = 14.628604639902
$ws.Cells.Item(19, 3).Value = 10.45762186984085
$ws.Cells.Item(19, 4).Value = 4.185135266915941
$ws.Cells.Item(19, 5).Value = 11.81637461187918
$ws.Cells.Item(19, 6).Value = 21.30749322614087
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 13).Value = 15.4950255334012
$ws.Cells.Item(19, 15).Value = 18.70441458932848
$ws.Cells.Item(20, 2).Value = 14.85996620816981
$ws.Cells.Item(20, 3).Value = 10.60944160363797
$ws.Cells.Item(20, 4).Value = 4.20201852502425
$ws.Cells.Item(20, 5).Value = 11.85716862668394
$ws.Cells.Item(20, 6).Value = 21.32820745896986
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 13).Value = 15.62262654735224
$ws.Cells.Item(20, 15).Value = 18.6911029126054
$ws.Cells.Item(21, 2).Value = 15.61207153456377
$ws.Cells.Item(21, 3).Value = 11.10286235973171
$ws.Cells.Item(21, 4).Value = 4.258100420375614
$ws.Cells.Item(21, 5).Value = 11.9992874760085
$ws.Cells.Item(21, 6).Value = 21.40839638792806
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 13).Value = 16.04687040177872
$ws.Cells.Item(21, 15).Value = 18.65490758399702
$ws.Cells.Item(22, 2).Value = 16.08448177233268
$ws.Cells.Item(22, 3).Value = 11.4127226860645
$ws.Cells.Item(22, 4).Value = 4.294239933918202
$ws.Cells.Item(22, 5).Value = 12.09579638627971
$ws.Cells.Item(22, 6).Value = 21.46865414012586
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 13).Value = 16.32033522414057
$ws.Cells.Item(22, 15).Value = 18.63764988466402
$ws.Cells.Item(23, 2).Value = 15.83411703479728
$ws.Cells.Item(23, 3).Value = 11.24851018139974
$ws.Cells.Item(23, 4).Value = 4.275000683146545
$ws.Cells.Item(23, 5).Value = 12.04396540770479
$ws.Cells.Item(23, 6).Value = 21.43577851316295
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 13).Value = 16.17475873923043
$ws.Cells.Item(23, 15).Value = 18.6462727686422
$ws.Cells.Item(24, 2).Value = 14.84735496661258
$ws.Cells.Item(24, 3).Value = 10.60116657250653
$ws.Cells.Item(24, 4).Value = 4.201093711234425
$ws.Cells.Item(24, 5).Value = 11.85490903931279
$ws.Cells.Item(24, 6).Value = 21.32702998952216
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 13).Value = 15.61563492792558
$ws.Cells.Item(24, 15).Value = 18.69180194461336
$ws.Cells.Item(25, 2).Value = 13.70310366021947
$ws.Cells.Item(25, 3).Value = 9.850071848488247
$ws.Cells.Item(25, 4).Value = 4.119389504970804
$ws.Cells.Item(25, 5).Value = 11.66744963801123
$ws.Cells.Item(25, 6).Value = 21.24368185889139
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 13).Value = 14.99929647513908
$ws.Cells.Item(25, 15).Value = 18.76810836994687
